$wb = $excel.ActiveWorkbook

# The "Users" sheet (sheet2) contains the Standard User's name in cell A2.
$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Ayati Arvind"

# Update the active selection on the Users sheet to D2 (matches target sheetView selection).
$usersSheet.Activate()
$usersSheet.Range("D2").Select()
